# Apply updated crypto price/volume figures (refreshed data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new display text. Column D ("Price") holds values that look
# numeric (e.g. "1.000", "0.9999") but must stay literal text, matching the
# workbook's existing inline-string formatting; Column E ("Volume(1h)") is
# always text (percentages padded with spaces).
$updates = @(
    @{ Cell = "D2"; Value = '25.916.77' }
    @{ Cell = "E2"; Value = '  -0.15%  ' }
    @{ Cell = "D3"; Value = '1.743.04' }
    @{ Cell = "E3"; Value = '  -0.84%  ' }
    @{ Cell = "D4"; Value = '1.000' }
    @{ Cell = "E4"; Value = '  -0.15%  ' }
    @{ Cell = "D5"; Value = '230.38' }
    @{ Cell = "E5"; Value = '  -2.21%  ' }
    @{ Cell = "E6"; Value = '  -0.07%  ' }
    @{ Cell = "D7"; Value = '0.5244' }
    @{ Cell = "E7"; Value = '  +0.65%  ' }
    @{ Cell = "D8"; Value = '0.2749' }
    @{ Cell = "E8"; Value = '  +0.38%  ' }
    @{ Cell = "D9"; Value = '39.47' }
    @{ Cell = "E9"; Value = '  -2.63%  ' }
    @{ Cell = "D10"; Value = '0.06151' }
    @{ Cell = "E10"; Value = '  -0.24%  ' }
    @{ Cell = "D11"; Value = '1.742.01' }
    @{ Cell = "E11"; Value = '  -0.95%  ' }
    @{ Cell = "D12"; Value = '0.07101' }
    @{ Cell = "E12"; Value = '  +0.87%  ' }
    @{ Cell = "E13"; Value = '  -2.36%  ' }
    @{ Cell = "D14"; Value = '0.6431' }
    @{ Cell = "E14"; Value = '  +1.13%  ' }
    @{ Cell = "D15"; Value = '4.530' }
    @{ Cell = "E15"; Value = '  +0.13%  ' }
    @{ Cell = "D16"; Value = '77.51' }
    @{ Cell = "E16"; Value = '  -0.28%  ' }
    @{ Cell = "D17"; Value = '1.000' }
    @{ Cell = "E17"; Value = '  -0.12%  ' }
    @{ Cell = "D18"; Value = '0.9999' }
    @{ Cell = "E18"; Value = '  -0.18%  ' }
    @{ Cell = "D19"; Value = '25.899.99' }
    @{ Cell = "E19"; Value = '  -0.22%  ' }
    @{ Cell = "E20"; Value = '  -0.60%  ' }
    @{ Cell = "D21"; Value = '0.000006674' }
    @{ Cell = "D22"; Value = '1.962.77' }
    @{ Cell = "E22"; Value = '  -1.66%  ' }
    @{ Cell = "D23"; Value = '4.300' }
    @{ Cell = "E23"; Value = '  +6.06%  ' }
    @{ Cell = "D24"; Value = '8.772' }
    @{ Cell = "E24"; Value = '  +3.48%  ' }
    @{ Cell = "D25"; Value = '5.165' }
    @{ Cell = "E25"; Value = '  -0.11%  ' }
    @{ Cell = "D26"; Value = '140.37' }
    @{ Cell = "E26"; Value = '  +1.02%  ' }
    @{ Cell = "D27"; Value = '1.520' }
    @{ Cell = "E27"; Value = '  +1.14%  ' }
    @{ Cell = "D28"; Value = '15.21' }
    @{ Cell = "E28"; Value = '  +0.66%  ' }
    @{ Cell = "D29"; Value = '1.793' }
    @{ Cell = "E29"; Value = '  -2.57%  ' }
    @{ Cell = "D30"; Value = '102.97' }
    @{ Cell = "E30"; Value = '  -0.10%  ' }
    @{ Cell = "D31"; Value = '0.08328' }
    @{ Cell = "E31"; Value = '  -0.58%  ' }
    @{ Cell = "D32"; Value = '3.725' }
    @{ Cell = "E32"; Value = '  +1.65%  ' }
    @{ Cell = "D33"; Value = '3.550' }
    @{ Cell = "E33"; Value = '  +3.61%  ' }
    @{ Cell = "D34"; Value = '0.04532' }
    @{ Cell = "E34"; Value = '  +1.77%  ' }
    @{ Cell = "D35"; Value = '2.618' }
    @{ Cell = "E35"; Value = '  -0.25%  ' }
    @{ Cell = "D36"; Value = '0.9770' }
    @{ Cell = "E36"; Value = '  -1.72%  ' }
    @{ Cell = "D37"; Value = '0.6192' }
    @{ Cell = "E37"; Value = '  +2.56%  ' }
    @{ Cell = "D38"; Value = '2.691' }
    @{ Cell = "E38"; Value = '  -1.27%  ' }
    @{ Cell = "D39"; Value = '0.01589' }
    @{ Cell = "E39"; Value = '  +0.33%  ' }
    @{ Cell = "D40"; Value = '1.924' }
    @{ Cell = "E40"; Value = '  -1.40%  ' }
    @{ Cell = "D41"; Value = '1.000' }
    @{ Cell = "E41"; Value = '  -0.09%  ' }
    @{ Cell = "D42"; Value = '100.25' }
    @{ Cell = "E42"; Value = '  -2.18%  ' }
    @{ Cell = "D43"; Value = '0.3868' }
    @{ Cell = "E43"; Value = '  +0.22%  ' }
    @{ Cell = "D44"; Value = '0.7327' }
    @{ Cell = "E44"; Value = '  -0.99%  ' }
    @{ Cell = "D45"; Value = '5.010' }
    @{ Cell = "E45"; Value = '  +2.10%  ' }
    @{ Cell = "D46"; Value = '0.05339' }
    @{ Cell = "E46"; Value = '  -3.09%  ' }
    @{ Cell = "D47"; Value = '0.1124' }
    @{ Cell = "E47"; Value = '  +1.18%  ' }
    @{ Cell = "D48"; Value = '6.263' }
    @{ Cell = "E48"; Value = '  -0.45%  ' }
    @{ Cell = "D49"; Value = '53.53' }
    @{ Cell = "E49"; Value = '  +2.27%  ' }
    @{ Cell = "D50"; Value = '30.17' }
    @{ Cell = "E50"; Value = '  +0.32%  ' }
    @{ Cell = "D51"; Value = '7.695' }
    @{ Cell = "E51"; Value = '  +3.87%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Cell.StartsWith("D")) {
        # Leading apostrophe forces text interpretation so price strings like
        # "1.000" / "0.9999" are not coerced into numbers; resetting the style
        # back to Normal strips the quote-prefix marker Excel adds, leaving a
        # plain General-formatted text cell (matching the original).
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
